# Updated symbol list on Mon Jan  9 08:58:17 UTC 2023 with GitHub Actions
#
# Refreshes the cached "Price" (column D) and "Volume(1h)" (column E)
# quotes for each coin row on Sheet1. Values are plain text in the
# workbook (not numbers/percentages), so they are written with a
# leading apostrophe to force Excel to keep them as literal text
# instead of auto-converting "277.54" to a number or "6.36%" to a
# percentage value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'277.54"
$ws.Range("E2").Value  = "'6.36%"

$ws.Range("D3").Value  = "'27.28"
$ws.Range("E3").Value  = "'1.49%"

$ws.Range("D4").Value  = "'4.789"
$ws.Range("E4").Value  = "'1.39%"

$ws.Range("D5").Value  = "'0.06280"
$ws.Range("E5").Value  = "'1.02%"

$ws.Range("E6").Value  = "'2.92%"

$ws.Range("D7").Value  = "'0.8812"
$ws.Range("E7").Value  = "'3.78%"

$ws.Range("D8").Value  = "'0.9464"
$ws.Range("E8").Value  = "'3.55%"

$ws.Range("D9").Value  = "'0.1458"
$ws.Range("E9").Value  = "'4.03%"

$ws.Range("D10").Value = "'0.05348"
$ws.Range("E10").Value = "'8.58%"

$ws.Range("D11").Value = "'0.07282"
$ws.Range("E11").Value = "'2.69%"

$ws.Range("D12").Value = "'0.03100"
$ws.Range("E12").Value = "'-0.39%"

$ws.Range("D13").Value = "'0.09074"
$ws.Range("E13").Value = "'0.30%"

$ws.Range("D14").Value = "'0.001557"
$ws.Range("E14").Value = "'1.40%"

$ws.Range("D15").Value = "'0.0006292"
$ws.Range("E15").Value = "'2.24%"

$ws.Range("D16").Value = "'0.005922"
$ws.Range("E16").Value = "'-0.80%"

$ws.Range("D17").Value = "'3.445"
$ws.Range("E17").Value = "'-0.12%"

$ws.Range("D18").Value = "'3.266"
$ws.Range("E18").Value = "'2.91%"

$ws.Range("E19").Value = "'5.32%"

$ws.Range("E20").Value = "'1.61%"

$ws.Range("E21").Value = "'0.12%"

$ws.Range("E22").Value = "'-6.17%"

$ws.Range("D23").Value = "'0.04328"
$ws.Range("E23").Value = "'1.73%"

$ws.Range("D24").Value = "'0.001182"
$ws.Range("E24").Value = "'-0.11%"

$ws.Range("E25").Value = "'5.22%"

$ws.Range("E26").Value = "'0.05%"

$ws.Range("D27").Value = "'0.0001692"
$ws.Range("E27").Value = "'3.16%"

$ws.Range("D40").Value = "'0.04055"
$ws.Range("E40").Value = "'3.03%"

$ws.Range("E41").Value = "'55.18%"

$ws.Range("D42").Value = "'0.1156"
$ws.Range("E42").Value = "'3.92%"

$ws.Range("D43").Value = "'0.002149"
$ws.Range("E43").Value = "'0.34%"

$ws.Range("D44").Value = "'0.01299"
$ws.Range("E44").Value = "'-2.41%"

$ws.Range("D45").Value = "'0.00005110"
$ws.Range("E45").Value = "'-1.01%"

$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("E46").Value = "'0.05%"

$ws.Range("E47").Value = "'851.60%"

$ws.Range("D49").Value = "'0.00002102"
$ws.Range("E49").Value = "'0.05%"

$ws.Range("D50").Value = "'0.0002002"
$ws.Range("E50").Value = "'0.05%"
